# edit.ps1
# Applies two changes described by the commit/diff:
#   1. "Added missing date to teacher meeting 3-19" -> the header's
#      "TEMP, 2015" placeholder text becomes "May 19, 2015" (typed as a
#      new run "May 19" immediately followed by the pre-existing ", 2015"
#      run, matching the way Word would keep the freshly typed text in
#      its own run).
#   2. Word's automatic "_GoBack" bookmark (marks the last editing
#      position before the file was saved) moves from where the previous
#      save had left it (end of the "Josh: Read in 1st and 2nd" line) to
#      the very start of the document, because the most recent edit
#      happened there.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Update the header date text: "TEMP, 2015" -> "May 19" + ", 2015"
# ---------------------------------------------------------------------
$headerRange = $d.Sections.Item(1).Headers.Item(1).Range

$headerFragment = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="007D0F06" w:rsidRPr="00404E5C" w:rsidRDefault="002E26A0"><w:pPr><w:rPr><w:color w:val="404040" w:themeColor="text1" w:themeTint="BF"/></w:rPr></w:pPr><w:r w:rsidRPr="00404E5C"><w:rPr><w:color w:val="404040" w:themeColor="text1" w:themeTint="BF"/></w:rPr><w:t>May 19</w:t></w:r><w:r><w:rPr><w:color w:val="404040" w:themeColor="text1" w:themeTint="BF"/></w:rPr><w:t>, 2015</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

[void]$headerRange.InsertXML($headerFragment)

# ---------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark to the start of the document
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$firstParagraphRange = $d.Paragraphs.Item(1).Range

$firstParagraphFragment = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="007D0F06" w:rsidRPr="00292367" w:rsidRDefault="002E26A0"><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r w:rsidRPr="00292367"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">Tune Squad Meeting Minutes </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

[void]$firstParagraphRange.InsertXML($firstParagraphFragment)

Write-Output "Header text now: $($headerRange.Text)"
Write-Output "GoBack bookmark present: $($d.Bookmarks.Exists('_GoBack'))"
